# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) -------------------------------------------------
# Copy the formatting of the existing "IP" header cell (H1, style index 1:
# bold font, thin border, centered) onto the two new header cells so they
# match the look of the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$ws.Range("H1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-37) ------------------------------------------------------
$data = @{
    2  = @(7, 9)
    3  = @(7, 7)
    4  = @(2, 5)
    5  = @(3, 5)
    6  = @(2, 3)
    7  = @(2, 5)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(8, 8)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(1, 6)
    15 = @(1, 4)
    16 = @(1, 5)
    17 = @(1, 4)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 6)
    21 = @(1, 6)
    22 = @(1, 6)
    23 = @(1, 5)
    24 = @(1, 5)
    25 = @(1, 5)
    26 = @(1, 6)
    27 = @(1, 6)
    28 = @(1, 6)
    29 = @(1, 7)
    30 = @(1, 5)
    31 = @(1, 4)
    32 = @(1, 5)
    33 = @(1, 6)
    34 = @(1, 5)
    35 = @(1, 5)
    36 = @(1, 3)
    37 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
